# Generate Report for Archive
# - Update status text from "Ready for handoff" to "In Translation"
# - Re-fit the "zh-cn"/"de-de"/"Status" columns that held that text

$wb = $excel.ActiveWorkbook

# Replace the old status text wherever it appears, being careful not to
# accidentally match boolean cells (PowerShell's -eq coerces a bare
# boolean left-hand side into matching any non-empty string), so only
# compare when the cell's value is actually a string.
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $v = $cell.Value2
        if (($v -is [string]) -and ($v -eq "Ready for handoff")) {
            $cell.Value = "In Translation"
        }
    }
}

# The columns that contained the status text need to be resized to fit
# the new (shorter) text.
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5   # column E ("zh-cn")
$overview.Columns.Item(6).ColumnWidth = 12.5   # column F ("de-de")

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5        # column C ("Status")

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5        # column C ("Status")
